# creating common data provider and adding log4j how logger
$wb = $excel.ActiveWorkbook

# --- rename the existing "AddCustomerTest" sheet to "ADDCUSTOMER" ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Name = "ADDCUSTOMER"

# --- add a new "OpenAccountTest" sheet right after it, as a common data provider ---
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccountTest"

# --- populate the new sheet with a small customer/currency table ---
$ws2.Range("A1").Value = "customer"
$ws2.Range("B1").Value = "currency"
$ws2.Range("A2").Value = "Anne Zimmermann"
$ws2.Range("B2").Value = "Real"

$ws2.Columns.Item(1).ColumnWidth = 17
$ws2.Range("B2").Select()

# --- go back to the first sheet, move the selection (was E3, now B6) and keep it active ---
$ws1.Activate()
$ws1.Range("B6").Select()
